# Add data from excel, ReadFromFile and YourInfoPage
# Populate the "user_details" sheet with a FirstName/LastName/PostalCode table,
# matching the formatting used on the "login" sheet (bordered cells, highlighted header row),
# and make "user_details" the active sheet/tab.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # login
$ws2 = $wb.Worksheets.Item(2)   # user_details

# --- Bring over the existing bordered-cell formatting used on the login sheet ---
$ws1.Range("A4").Copy() | Out-Null
$ws2.Range("A1:C2").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# --- Highlight the header row in yellow ---
$ws2.Range("A1:C1").Interior.Color = 65535

# --- Column widths for the new table ---
$ws2.Columns.Item(1).ColumnWidth = 11.5
$ws2.Columns.Item(2).ColumnWidth = 10.83
$ws2.Columns.Item(3).ColumnWidth = 14.67

# --- Fill in the values (column by column, header then data) ---
$ws2.Range("A1").Value = "FirstName "
$ws2.Range("B1").Value = "LastName"
$ws2.Range("A2").Value = "Neon"
$ws2.Range("B2").Value = "Test"
$ws2.Range("C1").Value = "PostalCode"
$ws2.Range("C2").Value = "'2000"

# --- Make user_details the active sheet/tab with C3 selected ---
$ws2.Activate() | Out-Null
$ws2.Range("C3").Select() | Out-Null
